# Update column F (dSF) values for the rows whose mean calculation
# changed after the data repull, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -1
    7  = -5
    10 = -5
    14 = -1
    16 = -3
    24 = 0
    30 = -1
    33 = 2
    34 = -2
    35 = 1
    51 = -2
    55 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
